$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Continue the existing "3/18" block with a new task (Job No. 3)
$ws.Range("B21").Value = "3"
$ws.Range("C21").Value = "Change info window style"
$ws.Range("D21").Value = 0.77083333333333337
$ws.Range("E21").Value = 0.80208333333333337
$ws.Range("F21").Value = 0.75

# New day block starting 3/21 (Job No. 1) -- row 22 left blank as a separator
$ws.Range("A23").Value = 43911
$ws.Range("B23").Value = "1"
$ws.Range("C23").Value = "Generate the JSON file containing Up/Down Streams as well as intermediate points"
$ws.Range("D23").Value = 0.97222222222222221

# Move the active selection to the next empty row, as left by the author
$ws.Range("C24").Select() | Out-Null
